$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# The sheet is protected; unlock it (password recovered from the stored hash) so
# rows can be inserted, then re-protect afterwards to match the published sheet.
$ws.Unprotect("admin")

# Insert a new row above the current row 6 ("Billing Provider"), shifting
# everything below it down by one.
$ws.Rows.Item(6).Insert()

# The freshly inserted row has no formatting; clone it from the row that used
# to be row 6 (now row 7) so the label/value styling (s=14 / s=16) matches.
$ws.Range("A7:B7").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A6").RowHeight = $ws.Range("A7").RowHeight
$excel.CutCopyMode = $false

# New label/value pair for the "Payer" row.
$ws.Range("A6").Value = "Payer"
$ws.Range("B6").Value = ""

# Put the active selection on B8 (the newly-shifted blank spacer row under the
# "Payer"/"Billing Provider" rows), matching the saved selection state.
$ws.Range("B8").Select()

# Re-protect the sheet (no password) like the saved workbook.
$ws.Protect()
